$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4834976253075638
$ws.Range("C2").Value = 0.1023584793711052
$ws.Range("D2").Value = 0.6596294088646744
$ws.Range("E2").Value = 0.2690594927575418
$ws.Range("G2").Value = 0.002490110008369273
$ws.Range("I2").Value = 0.8656079116674604
$ws.Range("J2").Value = 0.1400508006420935
$ws.Range("K2").Value = 0.5661745446029727
$ws.Range("N2").Value = 1.834853718740516
$ws.Range("O2").Value = 4.527969484952138

$ws.Range("B3").Value = 0.4462465166549805
$ws.Range("C3").Value = 0.09482582727186184
$ws.Range("D3").Value = 0.6493677274241065
$ws.Range("E3").Value = 0.2640200455365971
$ws.Range("G3").Value = 0.002492927336547925
$ws.Range("I3").Value = 0.8708834076828431
$ws.Range("J3").Value = 0.1367072967796545
$ws.Range("K3").Value = 0.5230310453022184
$ws.Range("N3").Value = 1.853804167883548
$ws.Range("O3").Value = 4.54153830315002

$ws.Range("B4").Value = 0.4235252953399993
$ws.Range("C4").Value = 0.09023551380785477
$ws.Range("D4").Value = 0.6434013599369734
$ws.Range("E4").Value = 0.261068473494511
$ws.Range("G4").Value = 0.002494750105087492
$ws.Range("I4").Value = 0.8745279677889357
$ws.Range("J4").Value = 0.1347330309150578
$ws.Range("K4").Value = 0.4967213293764416
$ws.Range("N4").Value = 1.866026601336567
$ws.Range("O4").Value = 4.552102723375981

$ws.Range("B5").Value = 0.414304678977544
$ws.Range("C5").Value = 0.08837371357628854
$ws.Range("D5").Value = 0.6410541965485095
$ws.Range("E5").Value = 0.2599015852808648
$ws.Range("G5").Value = 0.002495516333598295
$ws.Range("I5").Value = 0.8761150876646084
$ws.Range("J5").Value = 0.1339482788271198
$ws.Range("K5").Value = 0.4860457926786239
$ws.Range("N5").Value = 1.871154873043434
$ws.Range("O5").Value = 4.556969349617219

$ws.Range("B6").Value = 0.4127759397009356
$ws.Range("C6").Value = 0.08806509536512408
$ws.Range("D6").Value = 0.6406695392431061
$ws.Range("E6").Value = 0.2597099935163172
$ws.Range("G6").Value = 0.002495644982821719
$ws.Range("I6").Value = 0.8763847840240686
$ws.Range("J6").Value = 0.1338191657853756
$ws.Range("K6").Value = 0.4842759135571271
$ws.Range("N6").Value = 1.872015328614316
$ws.Range("O6").Value = 4.557811365067352

$ws.Range("B7").Value = 0.4234007865156002
$ws.Range("C7").Value = 0.09021036926641557
$ws.Range("D7").Value = 0.6433693642884748
$ws.Range("E7").Value = 0.2610525910334545
$ws.Range("G7").Value = 0.002494760343751623
$ws.Range("I7").Value = 0.874548959531861
$ws.Range("J7").Value = 0.1347223673987727
$ws.Range("K7").Value = 0.4965771689258816
$ws.Range("N7").Value = 1.866095165684317
$ws.Range("O7").Value = 4.552166082749267

$ws.Range("B8").Value = 0.4706223485501937
$ws.Range("C8").Value = 0.09975403148358453
$ws.Range("D8").Value = 0.6560218459847249
$ws.Range("E8").Value = 0.2672922937929982
$ws.Range("G8").Value = 0.002491062180793977
$ws.Range("I8").Value = 0.8673427613957401
$ws.Range("J8").Value = 0.1388816321104613
$ws.Range("K8").Value = 0.5512614636577098
$ws.Range("N8").Value = 1.84126603686185
$ws.Range("O8").Value = 4.532184525792218

$ws.Range("B9").Value = 0.5644070385962436
$ws.Range("C9").Value = 0.1187438984480593
$ws.Range("D9").Value = 0.6834836573873702
$ws.Range("E9").Value = 0.2806601767082526
$ws.Range("G9").Value = 0.002484544126294461
$ws.Range("I9").Value = 0.8564285774414344
$ws.Range("J9").Value = 0.1476627705274822
$ws.Range("K9").Value = 0.65991451121306
$ws.Range("N9").Value = 1.797232173833317
$ws.Range("O9").Value = 4.510723277789452

$ws.Range("B10").Value = 0.6340182944035178
$ws.Range("C10").Value = 0.1328633357142621
$ws.Range("D10").Value = 0.7052752205614752
$ws.Range("E10").Value = 0.2911727841118292
$ws.Range("G10").Value = 0.002480198275596091
$ws.Range("I10").Value = 0.8503725726031206
$ws.Range("J10").Value = 0.1544970507268459
$ws.Range("K10").Value = 0.7405936905808517
$ws.Range("N10").Value = 1.767718001386992
$ws.Range("O10").Value = 4.505770905354382

$ws.Range("B11").Value = 0.6658374460798484
$ws.Range("C11").Value = 0.139323134365327
$ws.Range("D11").Value = 0.7155396594751551
$ws.Range("E11").Value = 0.2961056877663211
$ws.Range("G11").Value = 0.002478316459113606
$ws.Range("I11").Value = 0.8480440656265031
$ws.Range("J11").Value = 0.1576897127366266
$ws.Range("K11").Value = 0.7774796230957293
$ws.Range("N11").Value = 1.754907196121973
$ws.Range("O11").Value = 4.505869229111511

$ws.Range("B12").Value = 0.6779081079864397
$ws.Range("C12").Value = 0.1417745580810674
$ws.Range("D12").Value = 0.7194770031974826
$ws.Range("E12").Value = 0.2979953135398361
$ws.Range("G12").Value = 0.002477617471593098
$ws.Range("I12").Value = 0.8472236621438896
$ws.Range("J12").Value = 0.1589107464732535
$ws.Range("K12").Value = 0.7914735396435617
$ws.Range("N12").Value = 1.750144594112852
$ws.Range("O12").Value = 4.506244708368655

$ws.Range("B13").Value = 0.67530752955156
$ws.Range("C13").Value = 0.1412463679965015
$ws.Range("D13").Value = 0.7186267847716863
$ws.Range("E13").Value = 0.2975873866487149
$ws.Range("G13").Value = 0.002477767406476902
$ws.Range("I13").Value = 0.8473976215085415
$ws.Range("J13").Value = 0.1586472393804996
$ws.Range("K13").Value = 0.7884585507086399
$ws.Range("N13").Value = 1.751166364991967
$ws.Range("O13").Value = 4.506148796195106

$ws.Range("B14").Value = 0.6668300799362612
$ws.Range("C14").Value = 0.1395247097941024
$ws.Range("D14").Value = 0.7158625771230618
$ws.Range("E14").Value = 0.2962607148074952
$ws.Range("G14").Value = 0.002478258680329889
$ws.Range("I14").Value = 0.8479753407542674
$ws.Range("J14").Value = 0.1577899265562479
$ws.Range("K14").Value = 0.7786303926779397
$ws.Range("N14").Value = 1.754513598366522
$ws.Range("O14").Value = 4.505893340194888

$ws.Range("B15").Value = 0.66164017422048
$ws.Range("C15").Value = 0.1384708250236031
$ws.Range("D15").Value = 0.7141759837848554
$ws.Range("E15").Value = 0.2954509077601273
$ws.Range("G15").Value = 0.002478561371339717
$ws.Range("I15").Value = 0.8483372013315815
$ws.Range("J15").Value = 0.1572663660290061
$ws.Range("K15").Value = 0.7726137343338166
$ws.Range("N15").Value = 1.756575412766995
$ws.Range("O15").Value = 4.505780920129041

$ws.Range("B16").Value = 0.6319418659525695
$ws.Range("C16").Value = 0.1324419085691204
$ws.Range("D16").Value = 0.7046114768107543
$ws.Range("E16").Value = 0.2908534370439639
$ws.Range("G16").Value = 0.002480323165512889
$ws.Range("I16").Value = 0.8505333293123272
$ws.Range("J16").Value = 0.1542900876323472
$ws.Range("K16").Value = 0.7381867771389352
$ws.Range("N16").Value = 1.76856760719429
$ws.Range("O16").Value = 4.505811804136897

$ws.Range("B17").Value = 0.6137616480638144
$ws.Range("C17").Value = 0.128752750110948
$ws.Range("D17").Value = 0.6988338790330033
$ws.Range("E17").Value = 0.288071610205435
$ws.Range("G17").Value = 0.002481428288737542
$ws.Range("I17").Value = 0.8519898192275264
$ws.Range("J17").Value = 0.1524856811818296
$ws.Range("K17").Value = 0.7171138518996258
$ws.Range("N17").Value = 1.77608209406208
$ws.Range("O17").Value = 4.506433049819094

$ws.Range("B18").Value = 0.6033192654546724
$ws.Range("C18").Value = 0.1266343102983001
$ws.Range("D18").Value = 0.6955438340667399
$ws.Range("E18").Value = 0.286485759765462
$ws.Range("G18").Value = 0.002482072884738008
$ws.Range("I18").Value = 0.8528676863058351
$ws.Range("J18").Value = 0.1514557136097068
$ws.Range("K18").Value = 0.7050106666194438
$ws.Range("N18").Value = 1.780462141441918
$ws.Range("O18").Value = 4.507011669220475

$ws.Range("B19").Value = 0.5997861425256872
$ws.Range("C19").Value = 0.125917639981509
$ws.Range("D19").Value = 0.6944355653305649
$ws.Range("E19").Value = 0.2859512546716232
$ws.Range("G19").Value = 0.002482292674314474
$ws.Range("I19").Value = 0.8531718087506093
$ws.Range("J19").Value = 0.1511083374399504
$ws.Range("K19").Value = 0.7009157412388447
$ws.Range("N19").Value = 1.781955094301775
$ws.Range("O19").Value = 4.507245583018715

$ws.Range("B20").Value = 0.6156954768646017
$ws.Range("C20").Value = 0.1291451089372515
$ws.Range("D20").Value = 0.6994454916138011
$ws.Range("E20").Value = 0.2883662727264067
$ws.Range("G20").Value = 0.002481309719826741
$ws.Range("I20").Value = 0.8518306195567433
$ws.Range("J20").Value = 0.1526769478228971
$ws.Range("K20").Value = 0.7193553033359024
$ws.Range("N20").Value = 1.775276169676376
$ws.Range("O20").Value = 4.506344012840856

$ws.Range("B21").Value = 0.6693195345790741
$ws.Range("C21").Value = 0.1400302608004722
$ws.Range("D21").Value = 0.7166731242320168
$ws.Range("E21").Value = 0.2966498033089806
$ws.Range("G21").Value = 0.00247811401237618
$ws.Range("I21").Value = 0.8478039852553678
$ws.Range("J21").Value = 0.1580414132685064
$ws.Range("K21").Value = 0.7815164590980714
$ws.Range("N21").Value = 1.75352803026178
$ws.Range("O21").Value = 4.505959192707252

$ws.Range("B22").Value = 0.704490686975987
$ws.Range("C22").Value = 0.1471748443090348
$ws.Range("D22").Value = 0.728226238320417
$ws.Range("E22").Value = 0.3021897059470859
$ws.Range("G22").Value = 0.002476104764143388
$ws.Range("I22").Value = 0.8455299679312418
$ws.Range("J22").Value = 0.1616175929069357
$ws.Range("K22").Value = 0.8222937558770411
$ws.Range("N22").Value = 1.73983067142769
$ws.Range("O22").Value = 4.507679319467968

$ws.Range("B23").Value = 0.6857079365566392
$ws.Range("C23").Value = 0.1433588716714667
$ws.Range("D23").Value = 0.7220332686514723
$ws.Range("E23").Value = 0.2992214218668536
$ws.Range("G23").Value = 0.002477169900071815
$ws.Range("I23").Value = 0.8467109194561928
$ws.Range("J23").Value = 0.159702494708597
$ws.Range("K23").Value = 0.8005164669438614
$ws.Range("N23").Value = 1.747093937933746
$ws.Range("O23").Value = 4.506580801251516

$ws.Range("B24").Value = 0.6148211627142928
$ws.Range("C24").Value = 0.1289677157009521
$ws.Range("D24").Value = 0.6991688832140994
$ws.Range("E24").Value = 0.2882330138772602
$ws.Range("G24").Value = 0.002481363296090669
$ws.Range("I24").Value = 0.8519024675364903
$ws.Range("J24").Value = 0.1525904530952005
$ws.Range("K24").Value = 0.7183419059724372
$ws.Range("N24").Value = 1.775640341603038
$ws.Range("O24").Value = 4.506383576629986

$ws.Range("B25").Value = 0.5389105337553133
$ws.Range("C25").Value = 0.1135772130164696
$ws.Range("D25").Value = 0.6757708094632449
$ws.Range("E25").Value = 0.2769225122281043
$ws.Range("G25").Value = 0.002486229318819577
$ws.Range("I25").Value = 0.8590365949222303
$ws.Range("J25").Value = 0.1452201545532006
$ws.Range("K25").Value = 0.6303705183752299
$ws.Range("N25").Value = 1.808646152000257
$ws.Range("O25").Value = 4.514630664497872

